$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update imputed values in column B and D for specific rows to reflect
# the refreshed RandomForest algorithm output.
$ws.Range("D4").Value = -7.9261
$ws.Range("B8").Value = 5.734199999999999
$ws.Range("B10").Value = 5.7577
$ws.Range("B12").Value = 5.416099999999997
$ws.Range("D12").Value = -5.960199999999999
$ws.Range("D15").Value = -8.225399999999997
$ws.Range("D17").Value = -8.377399999999996
$ws.Range("B18").Value = 6.473499999999996
$ws.Range("D26").Value = -7.442200000000006
$ws.Range("D27").Value = -7.9574
$ws.Range("D28").Value = -8.016999999999996
$ws.Range("B37").Value = 8.749600000000006
$ws.Range("D37").Value = -8.009099999999997
$ws.Range("D47").Value = -7.692999999999997
$ws.Range("B55").Value = 5.879499999999998
$ws.Range("D65").Value = -7.857600000000003
$ws.Range("B68").Value = 4.894599999999996
$ws.Range("D73").Value = -8.148899999999996
$ws.Range("B77").Value = 8.975900000000003
$ws.Range("B78").Value = 9.513300000000001
$ws.Range("B81").Value = 5.502200000000005
$ws.Range("B82").Value = 5.717099999999999
$ws.Range("D84").Value = -8.204699999999999
$ws.Range("D85").Value = -8.859199999999998
$ws.Range("D93").Value = -6.911399999999993
$ws.Range("D95").Value = -7.502399999999999
$ws.Range("D98").Value = -7.3094
$ws.Range("D99").Value = -7.925800000000005
$ws.Range("D101").Value = -8.119799999999994
